# update T1C Pool / fix minor bug / add preload in save update button.
# Adds a new "YEAR" column (Q) populated with 2025 for each data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column Q
$ws.Range("Q1").Value = "YEAR"

# Fill YEAR = 2025 for rows 2-10 (the 9 data rows)
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 17).Value = 2025
}

# Match the new cell style used for the YEAR values: Comma [0]-based style,
# with the number format reset back to General (numFmtId 0).
$ws.Range("Q2:Q10").Style = "Comma [0]"
$ws.Range("Q2:Q10").NumberFormat = "General"

# Update the view: scroll so column I is the left-most visible column and
# select the freshly added YEAR range.
$ws.Application.ActiveWindow.ScrollColumn = 9
$ws.Range("Q2:Q10").Select()
